# Updated cryptos list on Sat Jul 27 07:56:37 UTC 2024 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for every coin row,
# and re-sort a handful of neighbouring rows that swapped rank (their
# Coin/Link/Price/Volume all move together).
#
# Leading "'" on numeric-looking Price strings forces Excel to keep them as
# text (matching the source data, which can contain thousand-separator dots
# like "3.275.15") instead of silently coercing to a Double and losing
# formatting (e.g. "1.00" -> 1, "0.0690" -> 0.069).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.151.95"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "3.275.15"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'587.56"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'185.89"
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "'0.417"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "3.845.21"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'28.73"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "68.169.44"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "3.278.33"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'5.88"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'13.66"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'383.33"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "'7.73"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'71.51"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "'9.93"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").Value = "'0.185"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'5.84"
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").Value = "'22.98"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("D36").Value = "'162.62"
$ws.Range("D37").Value = "'1.87"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'6.80"
$ws.Range("E39").Value = "  +4.88%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.79"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Value = "'4.62"
$ws.Range("E41").Value = "  +5.78%  "
$ws.Range("D42").Value = "'2.61"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'41.56"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'349.32"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'25.61"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.661.09"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0690"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "'32.14"
$ws.Range("E49").Value = "  +5.40%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.103"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +1.98%  "
